# Apply cryptos.xlsx symbol-list update (Sat Dec 31 2022 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.70"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.47"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.130"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.531"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.019"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8171"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8400"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06952"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02841"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09385"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001518"

$ws.Range("E15").Value = "14OneONE"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006140"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.505"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03160"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.754"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04726"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001251"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004267"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009702"

$ws.Range("E27").Value = "26NitroExNTX"

$ws.Range("E28").Value = "27UpBotsUBXTWorstin24h"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006222"

$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008305"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005298"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002128"

